$p = $ppt.ActivePresentation

# Slide 7's notes page: collapse the split runs into a single run of text.
# (Text is set via an intermediate value first so the engine treats it as a
# real content change rather than a no-op, since the final text already
# equals the concatenation of the existing runs.)
$notesTr = $p.Slides.Item(7).NotesPage.Shapes.Item(2).TextFrame.TextRange
$notesTr.Text = "temp"
$notesTr.Text = "This is a blank slide: does it have a footer?"

# Slide titles made of multiple runs ("Slide" " " "N") -> a single run
# with the full text "Slide N". Same two-step trick is used for the same
# reason (rendered text is otherwise unchanged).
$t2 = $p.Slides.Item(2).Shapes.Item(1).TextFrame.TextRange
$t2.Text = "temp"
$t2.Text = "Slide 1"

$t4 = $p.Slides.Item(4).Shapes.Item(1).TextFrame.TextRange
$t4.Text = "temp"
$t4.Text = "Slide 3"

$t5 = $p.Slides.Item(5).Shapes.Item(1).TextFrame.TextRange
$t5.Text = "temp"
$t5.Text = "Slide 4"

$t6 = $p.Slides.Item(6).Shapes.Item(1).TextFrame.TextRange
$t6.Text = "temp"
$t6.Text = "Slide 5"
